# "Changed CRSTb to CRESET"
#
# The reset-signal label "CRSTb" was renamed to "CRESET" throughout the
# document (the FT232H/UPDuino pin-mapping list entry "D7 - CRSTb" and the
# JP1 jumper description that mentions the signal twice).
#
# Use Find/Replace across the whole document body so every occurrence is
# caught, then repeat the same replacement in every header/footer story so
# the edit is applied uniformly no matter which story range the text lives
# in.  wdReplaceAll = 2, wdFindContinue = 1.

$d = $word.ActiveDocument

function Replace-Everywhere($range, [string]$old, [string]$new) {
    $range.Find.ClearFormatting()
    $range.Find.Execute($old, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $new, 2) | Out-Null
}

# Main document story.
Replace-Everywhere $d.Content "CRSTb" "CRESET"

# Headers / footers (every section, every header/footer kind) in case the
# label also appears there.
foreach ($sec in $d.Sections) {
    foreach ($idx in 1, 2, 3) {
        $hdr = $sec.Headers($idx)
        if ($hdr.Exists) {
            Replace-Everywhere $hdr.Range "CRSTb" "CRESET"
        }
        $ftr = $sec.Footers($idx)
        if ($ftr.Exists) {
            Replace-Everywhere $ftr.Range "CRSTb" "CRESET"
        }
    }
}

$d.Save()
